$d = $word.ActiveDocument

# --- Fix 1: footer "Page X of Y" - convert the NUMPAGES w:fldSimple field into a
#     "complex" field (begin/instrText/separate/result/end), matching the
#     PAGE field right before it in the same paragraph. We rebuild the whole
#     footer paragraph via InsertXML so the resulting run/field structure is
#     exact. -------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)

$flds = $ftr.Range.Fields
$numPagesField = $null
for ($i = 1; $i -le $flds.Count; $i++) {
    if ($flds.Item($i).Code.Text -match "NUMPAGES") {
        $numPagesField = $flds.Item($i)
    }
}

if ($numPagesField -ne $null) {
    $ftrPara = $ftr.Range.Paragraphs.Item(1)
    $ftrParaRange = $ftrPara.Range

    $footerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:pPr><w:pStyle w:val="Footer"/></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Page </w:t></w:r>' +
      '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
      '<w:r><w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText></w:r>' +
      '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
      '<w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r>' +
      '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
      '<w:r><w:t xml:space="preserve"> of </w:t></w:r>' +
      '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
      '<w:r><w:instrText xml:space="preserve"> NUMPAGES   \* MERGEFORMAT </w:instrText></w:r>' +
      '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
      '<w:r><w:rPr><w:noProof/></w:rPr><w:t>4</w:t></w:r>' +
      '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
      '</w:p>'

    $ftrParaRange.InsertXML($footerXml)
}

# --- Fix 2: first-page header - drop the stray empty run (<w:r><w:t/></w:r>)
#     left behind in the otherwise-empty header paragraph. ----------------
$hdr2 = $sec.Headers.Item(2)
if ($hdr2.Exists -and $hdr2.Range.Text.Trim() -eq "") {
    $hdr2Para = $hdr2.Range.Paragraphs.Item(1)
    $hdr2ParaRange = $hdr2Para.Range

    $headerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:pPr><w:pStyle w:val="Header"/></w:pPr>' +
      '</w:p>'

    $hdr2ParaRange.InsertXML($headerXml)
}
